$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting to the three new rows (98, 99, 100) before filling in
#     values, copying it from existing rows that already carry the desired
#     look (Calibri 11 "automatic" text, no fill/border/number format).
#     xlPasteFormats = -4122
$xlPasteFormats = -4122

# Row 98 -> look like row 96 (name/username/status columns)
$ws.Range("A96").Copy() | Out-Null
$ws.Range("A98").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C96").Copy() | Out-Null
$ws.Range("C98").PasteSpecial($xlPasteFormats) | Out-Null

# Row 99 -> name column like row 96, status column like row 97 (last row)
$ws.Range("A96").Copy() | Out-Null
$ws.Range("A99").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C97").Copy() | Out-Null
$ws.Range("C99").PasteSpecial($xlPasteFormats) | Out-Null

# Row 100 -> status column like row 96; name/username left unstyled (like row 97)
$ws.Range("C96").Copy() | Out-Null
$ws.Range("C100").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- Fill in the three new joiner profiles. Username (column B) is written
#     before the full name (column A) on each row so the shared-string table
#     grows in the same order as the source workbook.
$ws.Range("B98").Value = "saboo_aparna1"
$ws.Range("A98").Value = "Aparna Saboo"
$ws.Range("C98").Value = 0

$ws.Range("B99").Value = "yashok42"
$ws.Range("A99").Value = "Ashok Yadav"
$ws.Range("C99").Value = 0

$ws.Range("B100").Value = "kanhu_panda"
$ws.Range("A100").Value = "Kanhucharan Panda"
$ws.Range("C100").Value = 0

# --- Move the active selection as recorded in the workbook after the edit.
$ws.Range("J97").Select() | Out-Null
